$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(5, 6, 7, 8, 11, 12, 13, 17, 18, 19, 20, 21, 24, 25, 26, 27, 28, 29, 31, 32, 33)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 15).Value = "Other"
}

$ws.Cells.Item(34, 15).Value = "NA"
